$wb = $excel.ActiveWorkbook

# --- Fees sheet: rename "Wire transfer" comment to "Currency conversion or wire transfer" ---
# "Wire transfer" was a misnomer (see commit message) - update every cell that
# still carries the old comment text, wherever it appears in the "Comment" column.
$feesSheet = $wb.Worksheets.Item("Fees")
$used = $feesSheet.UsedRange
$rowCount = $used.Rows.Count
for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $feesSheet.Cells.Item($r, 2)
    if ($cell.Value2 -eq "Wire transfer") {
        $cell.Value2 = "Currency conversion or wire transfer"
    }
}

# The new label is longer, so the "Comment" column needs to be widened to fit
# (mirrors Excel's best-fit column autosizing for the new text).
$feesSheet.Columns.Item(2).ColumnWidth = 32.333333

# --- Foreign Currencies sheet: update row 7 values ---
$fx = $wb.Worksheets.Item("Foreign Currencies")
$fx.Range("B7").Value2 = 155
$fx.Range("G7").Value2 = -10.13

# --- ELSTER - Summary sheet: update C7 value ---
$elster = $wb.Worksheets.Item("ELSTER - Summary")
$elster.Range("C7").Value2 = 16.86
